# Atualizado por script em 02-12-2023 14:46
# Appends 3 new match rows (88, 89, 90) to the Thai League 1 2023-2024
# odds sheet, mirroring the format/style of the preceding data row (87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, in column order A..V (A/Indice is recomputed below).
$newRows = @(
    @{
        pais = "thailand"; torneio = "thai-league-1"; temporada = "2023-2024"
        data_partida = 45262.5
        home = "Port MTI FC"; home_ft_gols = 2
        away = "Pathum United"; away_ft_gols = 3
        home_opening_odds = 1.81; home_opening_data_hora = "26/11/2023 13:12"
        home_closing_odds = 2.02; home_closing_data_hora = "02/12/2023 11:55"
        draw_opening_odds = 3.85; draw_opening_data_hora = "26/11/2023 13:12"
        draw_closing_odds = 3.78; draw_closing_data_hora = "02/12/2023 11:55"
        away_opening_odds = 4.07; away_opening_data_hora = "26/11/2023 13:12"
        away_closing_odds = 3.49; away_closing_data_hora = "02/12/2023 11:55"
        url_partida = "https://www.betexplorer.com/football/thailand/thai-league-1/port-mti-fc-pathum-united/0Qqe3naB/"
    },
    @{
        pais = "thailand"; torneio = "thai-league-1"; temporada = "2023-2024"
        data_partida = 45262.54166666666
        home = "Police Tero"; home_ft_gols = 1
        away = "Trat FC"; away_ft_gols = 3
        home_opening_odds = 1.96; home_opening_data_hora = "27/11/2023 12:42"
        home_closing_odds = 1.9; home_closing_data_hora = "02/12/2023 12:52"
        draw_opening_odds = 3.68; draw_opening_data_hora = "27/11/2023 12:42"
        draw_closing_odds = 3.89; draw_closing_data_hora = "02/12/2023 12:52"
        away_opening_odds = 3.63; away_opening_data_hora = "27/11/2023 12:42"
        away_closing_odds = 3.79; away_closing_data_hora = "02/12/2023 12:52"
        url_partida = "https://www.betexplorer.com/football/thailand/thai-league-1/police-tero-trat-fc/lnx80pqU/"
    },
    @{
        pais = "thailand"; torneio = "thai-league-1"; temporada = "2023-2024"
        data_partida = 45262.58333333334
        home = "Nakhon Pathom"; home_ft_gols = 1
        away = "Chonburi"; away_ft_gols = 0
        home_opening_odds = 3.38; home_opening_data_hora = "25/11/2023 14:12"
        home_closing_odds = 3.35; home_closing_data_hora = "02/12/2023 13:55"
        draw_opening_odds = 3.7; draw_opening_data_hora = "25/11/2023 14:12"
        draw_closing_odds = 3.59; draw_closing_data_hora = "02/12/2023 13:55"
        away_opening_odds = 1.97; away_opening_data_hora = "25/11/2023 14:12"
        away_closing_odds = 2.14; away_closing_data_hora = "02/12/2023 13:55"
        url_partida = "https://www.betexplorer.com/football/thailand/thai-league-1/nakhon-pathom-chonburi/6Nn31QTN/"
    }
)

$lastRow = $ws.UsedRange.Rows.Count

foreach ($row in $newRows) {
    $targetRow = $lastRow + 1

    # Copy the formatting (styles/number formats) of the last data row down
    # onto the new row, then overwrite with plain values below.
    $ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
    $ws.Range("A" + $targetRow + ":V" + $targetRow).PasteSpecial(-4122)

    $ws.Cells.Item($targetRow, 1).Value = $targetRow - 1
    $ws.Cells.Item($targetRow, 2).Value = $row.pais
    $ws.Cells.Item($targetRow, 3).Value = $row.torneio
    $ws.Cells.Item($targetRow, 4).Value = $row.temporada
    $ws.Cells.Item($targetRow, 5).Value = $row.data_partida
    $ws.Cells.Item($targetRow, 6).Value = $row.home
    $ws.Cells.Item($targetRow, 7).Value = $row.home_ft_gols
    $ws.Cells.Item($targetRow, 8).Value = $row.away
    $ws.Cells.Item($targetRow, 9).Value = $row.away_ft_gols
    $ws.Cells.Item($targetRow, 10).Value = $row.home_opening_odds
    $ws.Cells.Item($targetRow, 11).Value = $row.home_opening_data_hora
    $ws.Cells.Item($targetRow, 12).Value = $row.home_closing_odds
    $ws.Cells.Item($targetRow, 13).Value = $row.home_closing_data_hora
    $ws.Cells.Item($targetRow, 14).Value = $row.draw_opening_odds
    $ws.Cells.Item($targetRow, 15).Value = $row.draw_opening_data_hora
    $ws.Cells.Item($targetRow, 16).Value = $row.draw_closing_odds
    $ws.Cells.Item($targetRow, 17).Value = $row.draw_closing_data_hora
    $ws.Cells.Item($targetRow, 18).Value = $row.away_opening_odds
    $ws.Cells.Item($targetRow, 19).Value = $row.away_opening_data_hora
    $ws.Cells.Item($targetRow, 20).Value = $row.away_closing_odds
    $ws.Cells.Item($targetRow, 21).Value = $row.away_closing_data_hora
    $ws.Cells.Item($targetRow, 22).Value = $row.url_partida

    $lastRow = $targetRow
}
